# Apply the "daily process" schedule changes on the ScheduleDetails sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ScheduleDetails")
$ws.Activate()

# Toggle the Decesion_Value column:
#  D2 (Daily)   : No  -> Yes
#  D4 (Monthly) : Yes -> No
$ws.Range("D2").Value = "Yes"
$ws.Range("D4").Value = "No"

# Replicate the saved selection/window state from the authored workbook.
$ws.Range("D9").Select()
